# "Generate Report for Handoff"
# The localization status report is regenerated: the zh-cn / de-de rows move
# from "In Translation" to "Ready for handoff" and the handoff timestamps are
# refreshed, both on the per-language sheets and on the roll-up "Overview"
# sheet. The Status/Date columns are also widened to fit the new values.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Target column width (character units) for the widened Status / HO-datetime
# columns. Excel's ColumnWidth setter snaps to whole pixels (1/6-character
# steps at the default Calibri-11 metrics used by this workbook), so the
# nearest representable width to the authored 17.2159881591797 is fed in here.
$newColWidth = 16.333333333333332

# --- zh-cn sheet: Status + Latest Handoff Datetime -------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("G2").Value = "2016-07-26 08:28:22"
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet: Status + Latest Handoff Datetime -------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("G2").Value = "2016-07-26 08:28:33"
$dede.Columns.Item(3).ColumnWidth = $newColWidth

# --- Overview sheet: per-language status + generate date -------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-07-26 08:28:33"
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth
